$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'compression knee'
    'hockey pads men'
    'black basketball'
    'leg pads baseball'
    'running tights men'
    'snowboarding pads'
    'compression leggings youth boys'
    'knee pads wrestling men'
    'workout pants for girls'
    'basketball team gear'
    'basketball apparel for girls'
    'boys basketball compression leggings'
    'compressions pants'
    'bjj tights men'
    'compression basketball pants'
    'mens baseball pants'
    'best compression leggings'
    'softball pants black'
    'sliding pants youth'
    'basketballs pants'
    'sport knee pads'
    'men s knee pads'
    'volleyball knee pads black'
    'pack of leggings'
    'cycling knee protector'
    'knee pads sleeve basketball'
    'compression tights recovery'
    'volleyball pads'
    'padded tights for basketball'
    'small knee pads'
    'soccer apparel youth'
    'cheap knee pads'
    'baseball gear boys'
    'knee pads six six one'
    'compression pad'
    'free volleyball'
    'gym pad'
    'compression pants for boys basketball'
    'girl knee pads'
    'boys leggings sports youth'
    'youth football pants with pads small'
    'mens wrestling knee pads'
    'sport compression pants boys'
    'softball pants for girls'
    'compression leggings'
    'mens lacrosse pads'
    'softball pants youth girls black'
    'sport pants men'
    'thigh compression tights'
    'youth knee compression'
    'basketball sleeve knee pads'
    'leggings boys basketball'
    'basketball gear'
    'lacrosse youth pads'
    'running tights for men'
    'compressions leggings for men'
    'knee gel pads'
    'mens big and tall pants'
    'basketball protection'
    'knee braces for men xxl'
    'leggings with mesh girls'
    'little black pants guaranteed to fit'
    '3/4 pants men'
    'mens capri tights'
    'boys basketball'
    'girls volleyball spandex'
    'girls softball pants'
    'lcl knee support'
    'arthritis equipment'
    'volleyball knee pads'
    'running compression men'
    'baseball pants youth xl'
    'youth xs knee pads basketball'
    'under pants for basketball'
    'guy legging'
    'under armor warm pants'
    'spandex leggings men thermal'
    'mems thermal leggings'
    'men compression 3 4 pants'
    'men compression leggings 3 4'
    'men compression pants adidas'
    'men compression tights 3 4'
    'mens compression 3 4 tights'
    'mens compression 3 4 leggings'
    'eastbay baseball pants'
    'eastbay leggings'
    'eastbay tights'
    'mcdavid knee pads youth'
    'basketball leggings nike'
    'nba basketball pants'
    'basketball 3 4 leggings'
    'elbow knee wrist pads for youth'
    'telsa thermals men'
    'wintergear compression leggings men'
    'dry skin tights'
    'mens workoit tights'
    'track leggings'
    'track tights for men'
    'lotus leggings men'
    'basketball apparel youth'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
